$d = $word.ActiveDocument

# ------------------------------------------------------------
# 1) Straightforward text fixes / paragraph restructuring via
#    Find & Replace (whole-run anchors keep each match unique).
# ------------------------------------------------------------
$old = "is the only required undergraudate class I regularly teach, the others are electives and graduate courses. I typically have high enrollments in these elective classes (>20 students) and score well on my course evaluations. While the comments are ubiquitously supportive, I do get the occasional student comment that they don’t like the reverse classroom approach. These comments are difficult to address and are by far a minority. During COVID, I did have a few semesters with below (my) average evaluations. I attribute this to the social and interactive"
$new = "is the only required undergraduate class I regularly teach, the others are electives and graduate courses. I typically have high enrollments in these elective classes (>20 students) and score well on my course evaluations. While the comments are ubiquitously supportive, I do get the occasional student comment that they don’t like the reverse classroom approach. These comments are difficult to address and are by far a minority. During the COVID19 pandemic, I did have a few semesters with below (my) average evaluations. I attribute this to the social and interactive"
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Output "REPLACE FAILED: $old" }

$old = "part of these course being highly disadvantaged on Zoom. Since returning to in person classes post COVID, my course evaluation scores rebounded and even improved. I’ll also note that because of the larger enrollments and high response rates that all of my course evaluation scores represent good averages and are not artificially high due to only a few students reponding."
$new = "part of these course being highly disadvantaged on Zoom. Since returning to in person classes post COVID, my course evaluation scores rebounded and even improved."
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Output "REPLACE FAILED: $old" }

$old = "During the Summer 2023, I recieved academic development funds from the Cockrell School of Engineering"
$new = "During the Summer 2023, I received academic development funds from the Cockrell School of Engineering"
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Output "REPLACE FAILED: $old" }

$old = "which has shown continuous improvment with ratings of 3.7, 4.4, and 4.68 in three times I’ve taught the course since. This improvement was in spite of requiring the students to program for every assignment (21 assignments and 3 projects)!"
$new = "which has shown continuous improvement with ratings of 3.7, 4.4, and 4.68 in three times I’ve taught the course since. This improvement was in spite of requiring the students to program for every assignment (21 assignments and 3 projects)! I’ll also note that because of the larger enrollments and high response rates that all of my course evaluation scores represent good averages and are not artificially high due to only a few students responding."
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Output "REPLACE FAILED: $old" }

# ------------------------------------------------------------
# 2) Turn "animated computer visualization" into a hyperlink inside
#    the sentence "... Figure 2 shows an animated computer
#    visualization of tangent lines being drawn ...", splitting the
#    run into:
#      "shows an" | " " | [hyperlink]animated computer visualization[/hyperlink] | " " | "of tangent lines..."
# ------------------------------------------------------------
$anchor = $d.Content
$ok = $anchor.Find.Execute("shows an animated computer visualization of tangent lines", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { Write-Output "FIND FAILED: shows an animated computer visualization" }
$startPos = $anchor.Start

$p1 = $startPos + 8    # end of "shows an"
$p2 = $p1 + 1          # end of following space / start of "animated..."
$p3 = $p2 + 31         # end of "animated computer visualization"
$p4 = $p3 + 1          # end of following space / start of "of tangent..."

# Force a run-split at the "shows an" / " " boundary by toggling a
# formatting property on then off (leaves no visible formatting change).
$splitRange1 = $d.Range($p1, $p2)
$splitRange1.Bold = 1
$splitRange1.Bold = 0

# Force a run-split at the " " / "of tangent..." boundary the same way.
$splitRange2 = $d.Range($p3, $p4)
$splitRange2.Bold = 1
$splitRange2.Bold = 0

# Wrap "animated computer visualization" in a hyperlink.
$hlRange = $d.Range($p2, $p3)
if ($hlRange.Text -ne "animated computer visualization") {
    Write-Output "UNEXPECTED RANGE TEXT: $($hlRange.Text)"
}
$newLink = $d.Hyperlinks.Add($hlRange, "https://johnfoster.pge.utexas.edu/img/newton.gif", "", "", "animated computer visualization")
$hlRange.Style = "InternetLink"
